# Update the "想去人数" (F column) values on the "展览" and "全部类型"
# sheets to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 10
$ws1.Range("F4").Value = 13286
$ws1.Range("F7").Value = 680
$ws1.Range("F12").Value = 26
$ws1.Range("F13").Value = 13258
$ws1.Range("F15").Value = 580
$ws1.Range("F17").Value = 7938
$ws1.Range("F18").Value = 234
$ws1.Range("F28").Value = 196
$ws1.Range("F29").Value = 131
$ws1.Range("F30").Value = 360

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 10
$ws4.Range("F5").Value = 13286
$ws4.Range("F8").Value = 680
$ws4.Range("F13").Value = 26
$ws4.Range("F14").Value = 13258
$ws4.Range("F16").Value = 580
$ws4.Range("F18").Value = 7938
$ws4.Range("F19").Value = 234
$ws4.Range("F31").Value = 196
$ws4.Range("F32").Value = 131
$ws4.Range("F33").Value = 360
